$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.01253208636536152
$ws.Range("C2").Value = 2919.202174992006
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 14773364.14517103
$ws.Range("G2").Value = 14776302.07667549
